$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 4
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 2
$ws.Range("E8").Value = "30'"
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = "10'"
$excel.Calculate()
